$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value would be read back as a genuine number (e.g. "0.9998")
# are temporarily switched to a Text number format before the assignment so
# Excel keeps them as text, exactly like the other price cells in this sheet.
# The temporary format is removed again afterwards via ClearFormats so the
# cell style reverts back to the sheet default.
$textForcedRefs = @()

$ws.Range("D2").Value = '28.929.48'
$ws.Range("D3").Value = '1.832.79'
$ws.Range("E3").Value = '  -1.89%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9998'
$textForcedRefs += "D4"
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.85'
$textForcedRefs += "D5"
$ws.Range("E5").Value = '  +0.51%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6908'
$textForcedRefs += "D6"
$ws.Range("E6").Value = '  -1.78%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07691'
$textForcedRefs += "D8"
$ws.Range("E8").Value = '  -2.86%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3051'
$textForcedRefs += "D9"
$ws.Range("E9").Value = '  -2.59%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.45'
$textForcedRefs += "D10"
$ws.Range("E10").Value = '  -4.09%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07817'
$textForcedRefs += "D11"
$ws.Range("E11").Value = '  -0.23%  '
$ws.Range("D12").Value = '1.834.22'
$ws.Range("E12").Value = '  -1.29%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.078'
$textForcedRefs += "D13"
$ws.Range("E13").Value = '  -1.91%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '90.41'
$textForcedRefs += "D14"
$ws.Range("E14").Value = '  -3.72%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6795'
$textForcedRefs += "D15"
$ws.Range("E15").Value = '  -2.95%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.424'
$textForcedRefs += "D16"
$ws.Range("E16").Value = '  -1.62%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008311'
$textForcedRefs += "D17"
$ws.Range("E17").Value = '  -0.93%  '
$ws.Range("D18").Value = '28.938.30'
$ws.Range("E18").Value = '  -1.45%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '243.12'
$textForcedRefs += "D19"
$ws.Range("E19").Value = '  -3.88%  '
$ws.Range("D20").Value = '2.081.42'
$ws.Range("E20").Value = '  -1.50%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.73'
$textForcedRefs += "D21"
$ws.Range("E21").Value = '  -2.88%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9996'
$textForcedRefs += "D22"
$ws.Range("E22").Value = '  -0.04%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.468'
$textForcedRefs += "D23"
$ws.Range("E23").Value = '  -2.35%  '
$ws.Range("E24").Value = '  -0.07%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '163.05'
$textForcedRefs += "D25"
$ws.Range("E25").Value = '  +0.90%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1465'
$textForcedRefs += "D26"
$ws.Range("E26").Value = '  -5.61%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.794'
$textForcedRefs += "D27"
$ws.Range("E27").Value = '  -2.32%  '
$ws.Range("E28").Value = '  -3.52%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.547'
$textForcedRefs += "D29"
$ws.Range("E29").Value = '  +3.08%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.211'
$textForcedRefs += "D30"
$ws.Range("E30").Value = '  -2.39%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.150'
$textForcedRefs += "D31"
$ws.Range("E31").Value = '  -2.27%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.183'
$textForcedRefs += "D32"
$ws.Range("E32").Value = '  -2.64%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05105'
$textForcedRefs += "D33"
$ws.Range("E33").Value = '  -3.21%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7667'
$textForcedRefs += "D34"
$ws.Range("E34").Value = '  +2.44%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.837'
$textForcedRefs += "D35"
$ws.Range("E36").Value = '  -2.61%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.685'
$textForcedRefs += "D37"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01849'
$textForcedRefs += "D38"
$ws.Range("E38").Value = '  -1.54%  '
$ws.Range("D39").Value = '1.232.95'
$ws.Range("E39").Value = '  -3.05%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.698'
$textForcedRefs += "D40"
$ws.Range("E40").Value = '  -2.59%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9379'
$textForcedRefs += "D41"
$ws.Range("E41").Value = '  +5.12%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '107.71'
$textForcedRefs += "D42"
$ws.Range("E42").Value = '  -0.75%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9994'
$textForcedRefs += "D43"
$ws.Range("E43").Value = '  -0.14%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.660'
$textForcedRefs += "D44"
$ws.Range("E44").Value = '  -5.64%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '9.565'
$textForcedRefs += "D45"
$ws.Range("E45").Value = '  -0.58%  '
$ws.Range("E46").Value = '  -4.20%  '
$ws.Range("D47").Value = '1.981.22'
$ws.Range("E47").Value = '  -1.60%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5172'
$textForcedRefs += "D48"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '64.41'
$textForcedRefs += "D49"
$ws.Range("E49").Value = '  -9.28%  '
$ws.Range("E50").Value = '  -2.70%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4193'
$textForcedRefs += "D51"
$ws.Range("E51").Value = '  -2.44%  '

foreach ($ref in $textForcedRefs) {
    $ws.Range($ref).ClearFormats()
}
